# Auto-generated edit script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.248.70"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "2.367.64"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.04"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  -0.54%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.15"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  -5.60%  "
$__style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.07%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  -3.11%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.86"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  -5.27%  "
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0919"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  -2.52%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.47"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("E13").Value = "  -0.10%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.979"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  -4.66%  "
$ws.Range("D15").Value = "2.727.44"
$ws.Range("E15").Value = "  -2.12%  "
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.40"
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = "  -3.88%  "
$ws.Range("D17").Value = "2.371.90"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").Value = "45.192.21"
$ws.Range("E18").Value = "  -1.29%  "
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.49"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = "  +15.45%  "
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.28"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  -4.83%  "
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("E22").Value = "  +1.48%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.16"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  -2.73%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.70"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  -1.75%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -1.28%  "
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.49"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("E29").Value = "  -1.59%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.37"
$ws.Range("D30").Style = $__style
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0947"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  -2.06%  "
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.23"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  -5.17%  "
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.90"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  -5.23%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.71"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  -5.85%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.09"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  +6.64%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.03"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("E41").Value = "  -3.79%  "
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.62"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  -2.44%  "
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.10"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.874.01"
$ws.Range("E44").Value = "  +13.07%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.96"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("E46").Value = "  -5.42%  "
$ws.Range("E47").Value = "  +0.15%  "
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.96"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +1.82%  "
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.03"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  +6.10%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.99"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  -5.65%  "
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.17"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  -2.72%  "
